$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells currently holding "0.6632 0.0184" -> change to "0.66 0.0184"
$ws.Range("G4:I6").Value = "0.66 0.0184"

# Cells currently holding "0.3364 -0.0184" -> change to "0.3398 -0.0184"
$ws.Range("J7:O9").Value = "0.3398 -0.0184"
